$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "channel column in write module added": insert a new first column (A)
# for the channel name, shifting the existing N / Inversion / Numero de
# combinaciones / Cobertura columns from A:D to B:E.
$ws.Columns("A:A").Insert()

# Header row
$ws.Range("A1").Value = "Channel"
$ws.Range("B1").Value = "N"
$ws.Range("C1").Value = "Inversión"
$ws.Range("D1").Value = "Número de combinaciones"
$ws.Range("E1").Value = "Cobertura"

# Row 2 - Televisión
$ws.Range("A2").Value = "Televisión"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 543750
$ws.Range("D2").Value = 1216864
$ws.Range("E2").Value = 0.8605847883429888

# Row 3 - Digital Video
$ws.Range("A3").Value = "Digital Video"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 425000

# Row 4 - Cine
$ws.Range("A4").Value = "Cine"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 225000

# Row 5 - BVOD
$ws.Range("A5").Value = "BVOD"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 6250
